# Update the "fund_map" worksheet: restore the fundManager_raw/fundManager
# mapping table (clean label text, re-sorted/re-derived entries) and
# move the active selection to E26 to match the saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fund_map")

$ws.Range("A1").Value = "fundManager_raw"
$ws.Range("B1").Value = "fundManager"
$ws.Range("A2").Value = "Azzad Funds"
$ws.Range("B2").Value = "Azzad"
$ws.Range("A3").Value = "BARONFUNDS"
$ws.Range("B3").Value = "Baron"
$ws.Range("A4").Value = "BlackRock-Advised Funds"
$ws.Range("B4").Value = "BlackRock"
$ws.Range("A5").Value = "BlackRock-advised Funds"
$ws.Range("B5").Value = "BlackRock"
$ws.Range("A6").Value = "Brighthouse Funds Trust I"
$ws.Range("B6").Value = "Brighthouse"
$ws.Range("A7").Value = "Brighthouse Funds Trust II"
$ws.Range("B7").Value = "Brighthouse"
$ws.Range("A8").Value = "Brinker Capital Destinations Trust"
$ws.Range("B8").Value = "Brinker Capital"
$ws.Range("A9").Value = "CRMCFNDGRP"
$ws.Range("B9").Value = "Capital Group"
$ws.Range("A10").Value = "Clipper Funds Trust"
$ws.Range("B10").Value = "Clipper"
$ws.Range("A11").Value = "Davis Funds"
$ws.Range("B11").Value = "Davis "
$ws.Range("A12").Value = "Delaware Funds by Macquarie"
$ws.Range("B12").Value = "Delaware Funds"
$ws.Range("A13").Value = "FIRSTPACAD"
$ws.Range("B13").Value = "First Pacific Advisors"
$ws.Range("A14").Value = "Federated Hermes Funds"
$ws.Range("B14").Value = "Federated Hermes"
$ws.Range("A15").Value = "Fidelity Group of Funds"
$ws.Range("B15").Value = "Fidelity"
$ws.Range("A16").Value = "Franklin Templeton Group of Funds"
$ws.Range("B16").Value = "Franklin Templeton"
$ws.Range("A17").Value = "GREAT-WEST FUNDS INC"
$ws.Range("B17").Value = "Great-West"
$ws.Range("A18").Value = "Goldman Sachs Fund Complex"
$ws.Range("B18").Value = "Goldman Sachs"
$ws.Range("A19").Value = "Guggenheim"
$ws.Range("B19").Value = "Guggenheim"
$ws.Range("A20").Value = "HARTFORD FUNDS"
$ws.Range("B20").Value = "Hartford Funds"
$ws.Range("A21").Value = "INVESCOFDS"
$ws.Range("B21").Value = "Invesco"
$ws.Range("A22").Value = "JNL Fund Complex"
$ws.Range("B22").Value = "Jackson National"
$ws.Range("A23").Value = "John Hancock Group of Funds"
$ws.Range("B23").Value = "John Hancock"
$ws.Range("A24").Value = "LINCOLNTRS"
$ws.Range("B24").Value = "Lincoln Investment"
$ws.Range("A25").Value = "MASSMUTUAL FUNDS"
$ws.Range("B25").Value = "MassMutual"
$ws.Range("A26").Value = "Neuberger Berman"
$ws.Range("B26").Value = "Neuberger Berman"
$ws.Range("A27").Value = "PENN SERIES FUNDS INC"
$ws.Range("B27").Value = "Penn Mutual"
$ws.Range("A28").Value = "Principal Funds"
$ws.Range("B28").Value = "Principal Funds"
$ws.Range("A29").Value = "Publicly registered funds managed by Blackstone Alternative Asset Management L.P. Blackstone Alternative Investment Advisors LLC"
$ws.Range("B29").Value = "Blackstone"
$ws.Range("A30").Value = "SUNAMERICA"
$ws.Range("B30").Value = "Sun America"
$ws.Range("A31").Value = "Selected Funds"
$ws.Range("B31").Value = "Selected Funds"
$ws.Range("A32").Value = "TROWEPRICE"
$ws.Range("B32").Value = "T. Rowe Price"
$ws.Range("A33").Value = "VALIC Company"
$ws.Range("B33").Value = "VALIC"
$ws.Range("A34").Value = "Voya mutual funds"
$ws.Range("B34").Value = "Voya"

[void]$ws.Range("E26").Select()
